$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "An image"
